$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# "Experimental" row (row 7): set the Value cell (B7) to the literal text "true".
# A direct $c.Value = "true" would be auto-detected as the Boolean TRUE, but the
# source workbook stores this as a plain shared string, so enter it as a formula
# that evaluates to the text "true" and then convert that formula to its static
# value via copy / paste-special (values only) - this keeps the cell's stored
# type as a string instead of a boolean.
$b7 = $ws1.Cells.Item(7, 2)
$b7.Formula = '="true"'
$b7.Copy()
$b7.PasteSpecial(-4163)

# "Date" row (row 8): update the recorded timestamp for this value set export.
$ws1.Range("B8").Value = "2023-02-01T09:05:11-06:00"
